# Generate Report for Handback
#
# - Status text changes from "Ready for handoff" to
#   "Handed back: in sync with en-US" (every cell that shows that status).
# - The zh-cn / de-de sheets each grow two new "handback" columns
#   (E = Latest Handback File, F = Latest Handback DateTime... i.e. the
#   new target-file / target-file-hyperlink columns) for the two real
#   localized-file rows, mirroring the existing Source File (A) and
#   Latest Handoff File (C) hyperlinks.
# - de-de's "Latest Handback DateTime" column (G) moves from the
#   "0001-01-01 00:00:00" placeholder to the real handback timestamp.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 1. Status column: "Ready for handoff" -> "Handed back: in sync with en-US"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B2").Value = $newStatus
$wsZh.Range("B3").Value = $newStatus

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B2").Value = $newStatus
$wsDe.Range("B3").Value = $newStatus

# ---------------------------------------------------------------------
# 2. zh-cn sheet: add Latest Handback File (E) / Latest Handback
#    DateTime-as-file-link (F) hyperlinks for rows 2 and 3.
# ---------------------------------------------------------------------
$mdName = "b47b2731-05f2-4a1b-b1d5-bc1b0a1b1e7b.md"
$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/928f039c3066b5232b5c30402e8df4e21684d692/e2e/b47b2731-05f2-4a1b-b1d5-bc1b0a1b1e7b.md"

$zhXlfName = "b47b2731-05f2-4a1b-b1d5-bc1b0a1b1e7b.faa9ee56caa15bab7b3c7ef02717e49de8166f3e.zh-cn.xlf"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1782345dd02847c54c474d779cfd31dee1f37bb9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b47b2731-05f2-4a1b-b1d5-bc1b0a1b1e7b.faa9ee56caa15bab7b3c7ef02717e49de8166f3e.zh-cn.xlf"

foreach ($row in 2, 3) {
    $eCell = $wsZh.Range("E$row")
    $eCell.Value = $mdName
    $wsZh.Hyperlinks.Add($eCell, $mdUrl, "", "", $mdName)

    $fCell = $wsZh.Range("F$row")
    $fCell.Value = $zhXlfName
    $wsZh.Hyperlinks.Add($fCell, $zhXlfUrl, "", "", $zhXlfName)
}

# ---------------------------------------------------------------------
# 3. de-de sheet: add Latest Handback File (E) / Latest Handback
#    DateTime-as-file-link (F) hyperlinks for rows 2 and 3, and update
#    the Latest Handback DateTime (G) to the real handback timestamp.
# ---------------------------------------------------------------------
$deXlfName = "b47b2731-05f2-4a1b-b1d5-bc1b0a1b1e7b.faa9ee56caa15bab7b3c7ef02717e49de8166f3e.de-de.xlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/75860fe60f22a8e07214073c4cc32b900c857fd8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b47b2731-05f2-4a1b-b1d5-bc1b0a1b1e7b.faa9ee56caa15bab7b3c7ef02717e49de8166f3e.de-de.xlf"

$handbackDateTime = "2016-03-09 12:59:05"

foreach ($row in 2, 3) {
    $eCell = $wsDe.Range("E$row")
    $eCell.Value = $mdName
    $wsDe.Hyperlinks.Add($eCell, $mdUrl, "", "", $mdName)

    $fCell = $wsDe.Range("F$row")
    $fCell.Value = $deXlfName
    $wsDe.Hyperlinks.Add($fCell, $deXlfUrl, "", "", $deXlfName)

    $wsDe.Range("G$row").Value = $handbackDateTime
}

Write-Host "Handback report generated."
